$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the special "closing" bottom-border formatting (currently on
# row 26, the last data row) by copying it onto row 20, which will become
# the new last data row once the old rows are removed. ---
$ws.Range("B26:J26").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# --- Remove the old employee rows (21-26). The blank gap rows (27-30) were
# never materialised in the sheet, so deleting just 21-26 naturally shifts
# the signature block (old rows 31-32) up to rows 25-26, matching the
# target layout. ---
$ws.Rows("21:26").Delete()

# --- Update the summary header values ---
$ws.Range("E11").Value2 = 320000
$ws.Range("C13").Value2 = 2
$ws.Range("F13").Value2 = 4

# --- Replace the worker table contents with the new data ---
$ws.Range("B16").Value2 = "PPT"
$ws.Range("C16").Value2 = "6243847"
$ws.Range("D16").Value2 = "RONALD JARABA HERRERA"
$ws.Range("E16").Value2 = "2507"
$ws.Range("F16").Value2 = 64000
$ws.Range("G16").Value2 = 1600000

$ws.Range("B17").Value2 = "PPT"
$ws.Range("C17").Value2 = "6243847"
$ws.Range("D17").Value2 = "RONALD JARABA HERRERA"
$ws.Range("E17").Value2 = "2506"
$ws.Range("F17").Value2 = 64000
$ws.Range("G17").Value2 = 1600000

$ws.Range("B18").Value2 = "PPT"
$ws.Range("C18").Value2 = "6243847"
$ws.Range("D18").Value2 = "RONALD JARABA HERRERA"
$ws.Range("E18").Value2 = "2505"
$ws.Range("F18").Value2 = 64000
$ws.Range("G18").Value2 = 1600000

$ws.Range("B19").Value2 = "PPT"
$ws.Range("C19").Value2 = "6243847"
$ws.Range("D19").Value2 = "RONALD JARABA HERRERA"
$ws.Range("E19").Value2 = "2504"
$ws.Range("F19").Value2 = 64000
$ws.Range("G19").Value2 = 1600000

$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1002199824"
$ws.Range("D20").Value2 = "MARIO MELENDEZ STEVENSON"
$ws.Range("E20").Value2 = "2505"
$ws.Range("F20").Value2 = 64000
$ws.Range("G20").Value2 = 1600000
